$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates from the latest cryptos data pull.
# Cells whose new text renders as a bare number (e.g. "26.70") are force-
# typed as Text so Excel keeps them as strings (preserving formats like
# trailing zeros) instead of silently coercing them to numeric values.

$ws.Range("D2").Value = '62.448.07'
$ws.Range("E2").Value = '  +4.07%  '
$ws.Range("D3").Value = '3.335.07'
$ws.Range("E3").Value = '  +4.18%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.24'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.33'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.335.83'
$ws.Range("E8").Value = '  +4.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '3.910.52'
$ws.Range("E13").Value = '  +4.21%  '
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("E16").Value = '  +2.76%  '
$ws.Range("D17").Value = '62.471.82'
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("D18").Value = '3.315.72'
$ws.Range("E18").Value = '  +2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.76'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.37'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.48'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.96%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.88'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").Value = '  +5.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.55%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '0.0₃0942'
$ws.Range("E29").Value = '  +4.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.59'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.91%  '
$ws.Range("E31").Value = '  +3.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.59'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.84'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("E34").Value = '  +7.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.71'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.14'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.46'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +8.39%  '
$ws.Range("E38").Value = '  +12.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.80'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0736'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.68%  '
$ws.Range("D41").Value = '2.791.93'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  +7.76%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.41'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.24'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.741'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.18%  '
$ws.Range("E46").Value = '  +4.58%  '
$ws.Range("D47").Value = '3.376.66'
$ws.Range("E47").Value = '  +4.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.90'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.18%  '
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '286.71'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.01%  '
